$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 10-16: existing rows, scheme reassigned + values updated ---
# Row 10: Gaussian-Quadrature
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("C10").Value = 1.431032689423501
$ws.Range("D10").Value = 0.9794719751908326
$ws.Range("E10").Value = 1.214231639073742
$ws.Range("F10").Value = 0.8408026335015744
$ws.Range("G10").Value = 1.431032689423501
$ws.Range("H10").Value = 0.9794719751908326
$ws.Range("I10").Value = 1.091868016818288
$ws.Range("J10").Value = 0.8489512606847928
$ws.Range("K10").Value = 1.019012699760681
$ws.Range("L10").Value = 0.8639407924363697
$ws.Range("M10").Value = 1.431047313487744
$ws.Range("N10").Value = 1.096851807132287
$ws.Range("O10").Value = 1.116384734297413
$ws.Range("P10").Value = 1.036163963361223

# Row 11: Spiral-90deg-10rot-5space
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value = 0.9905450009565365
$ws.Range("D11").Value = 0.04764770810301898
$ws.Range("E11").Value = 1.667243911576986
$ws.Range("F11").Value = 0.8722746534224871
$ws.Range("G11").Value = 0.9905450009565365
$ws.Range("H11").Value = 0.04764770810301898
$ws.Range("I11").Value = 1.316676374834493
$ws.Range("J11").Value = 1.036123905525455
$ws.Range("K11").Value = 1.044071865548856
$ws.Range("L11").Value = 0.398367679958977
$ws.Range("M11").Value = 0.9905450009565365
$ws.Range("N11").Value = 0.8574458098400026
$ws.Range("O11").Value = 0.8944278185147573
$ws.Range("P11").Value = 0.9216188874908513

# Row 12: Spiral-90deg-15rot-5space
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value = 0.9835438214198388
$ws.Range("D12").Value = 0.04775026955329389
$ws.Range("E12").Value = 1.665651853766884
$ws.Range("F12").Value = 0.8742351688860515
$ws.Range("G12").Value = 0.9835438214198388
$ws.Range("H12").Value = 0.04775026955329389
$ws.Range("I12").Value = 1.315570487665723
$ws.Range("J12").Value = 1.03819736992039
$ws.Range("K12").Value = 1.04404674607472
$ws.Range("L12").Value = 0.3992774766780214
$ws.Range("M12").Value = 0.9835438214198388
$ws.Range("N12").Value = 0.8567010616600887
$ws.Range("O12").Value = 0.892795278406517
$ws.Range("P12").Value = 0.9210341492456152

# Row 13: Spiral-90deg-10rot-3space
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value = 0.9889001708670686
$ws.Range("D13").Value = 0.04764538641199453
$ws.Range("E13").Value = 1.667342640681802
$ws.Range("F13").Value = 0.8726856983582185
$ws.Range("G13").Value = 0.9889001708670686
$ws.Range("H13").Value = 0.04764538641199453
$ws.Range("I13").Value = 1.316567295785394
$ws.Range("J13").Value = 1.036749052037742
$ws.Range("K13").Value = 1.043443090767422
$ws.Range("L13").Value = 0.3984106438106491
$ws.Range("M13").Value = 0.9889001708670686
$ws.Range("N13").Value = 0.8574940135468984
$ws.Range("O13").Value = 0.894143474079771
$ws.Range("P13").Value = 0.9214679973400363

# Row 14: NoRotation-tilt60deg
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("C14").Value = 0.320187999999999
$ws.Range("D14").Value = 0.03040000000000005
$ws.Range("E14").Value = 2.212935999999999
$ws.Range("F14").Value = 0.9324560000000014
$ws.Range("G14").Value = 0.320187999999999
$ws.Range("H14").Value = 0.03040000000000005
$ws.Range("I14").Value = 1.462464000000001
$ws.Range("J14").Value = 0.8357559999999996
$ws.Range("K14").Value = 1.356939999999997
$ws.Range("L14").Value = 0.3413400000000005
$ws.Range("M14").Value = 0.320187999999999
$ws.Range("N14").Value = 1.121668
$ws.Range("O14").Value = 0.873995
$ws.Range("P14").Value = 0.9365599999999998

# Row 15: Rotation-NoTilt
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("C15").Value = 0.43
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 2.442362500000002
$ws.Range("F15").Value = 0.8126625000000011
$ws.Range("G15").Value = 0.43
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 1.610862499999995
$ws.Range("J15").Value = 0.5079250000000002
$ws.Range("K15").Value = 1.689037499999998
$ws.Range("L15").Value = 0.18
$ws.Range("M15").Value = 0.43
$ws.Range("N15").Value = 1.221181250000001
$ws.Range("O15").Value = 0.9212562500000009
$ws.Range("P15").Value = 0.9591062499999996

# Row 16: Rotation-60detTilt
$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("C16").Value = 0.6974514122752004
$ws.Range("D16").Value = 0.4004946589696001
$ws.Range("E16").Value = 1.809522171289599
$ws.Range("F16").Value = 0.8839243951104017
$ws.Range("G16").Value = 0.6974514122752004
$ws.Range("H16").Value = 0.4004946589696001
$ws.Range("I16").Value = 1.345522969702401
$ws.Range("J16").Value = 0.7269277276160011
$ws.Range("K16").Value = 1.371488815513599
$ws.Range("L16").Value = 0.5222412793856009
$ws.Range("M16").Value = 0.6973998473215998
$ws.Range("N16").Value = 1.105008415129599
$ws.Range("O16").Value = 0.9478481594112002
$ws.Range("P16").Value = 0.9696966787328003

# --- Rows 17-19: brand new rows ---
# Row 17: HexGrid-90degTilt5degRes
$ws.Range("A17").Value = 15
$ws.Range("A17").Font.Bold = $true
$ws.Range("A17").HorizontalAlignment = -4108
$ws.Range("A17").VerticalAlignment = -4160
$ws.Range("A17").Borders.LineStyle = 1
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C17").Value = 0.9935434098458009
$ws.Range("D17").Value = 0.9981320030533807
$ws.Range("E17").Value = 0.9902981108190573
$ws.Range("F17").Value = 0.9870418116290753
$ws.Range("G17").Value = 0.9935434098458009
$ws.Range("H17").Value = 0.9981320030533807
$ws.Range("I17").Value = 0.9921029124657249
$ws.Range("J17").Value = 0.9907313663900262
$ws.Range("K17").Value = 0.9896050237008106
$ws.Range("L17").Value = 0.9907160443981665
$ws.Range("M17").Value = 0.9935154270865922
$ws.Range("N17").Value = 0.9942150569362189
$ws.Range("O17").Value = 0.9922538338368285
$ws.Range("P17").Value = 0.9915213352877553

# Row 18: HexGrid-90degTilt22p5degRes
$ws.Range("A18").Value = 16
$ws.Range("A18").Font.Bold = $true
$ws.Range("A18").HorizontalAlignment = -4108
$ws.Range("A18").VerticalAlignment = -4160
$ws.Range("A18").Borders.LineStyle = 1
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18").Value = 0.9536572955620595
$ws.Range("D18").Value = 0.9819509516483769
$ws.Range("E18").Value = 1.011086090386303
$ws.Range("F18").Value = 1.063865736471009
$ws.Range("G18").Value = 0.9536572955620595
$ws.Range("H18").Value = 0.9819509516483769
$ws.Range("I18").Value = 0.9467545999943179
$ws.Range("J18").Value = 0.9909727471333972
$ws.Range("K18").Value = 1.015379311671147
$ws.Range("L18").Value = 1.007800727648787
$ws.Range("M18").Value = 0.9536572955620595
$ws.Range("N18").Value = 0.9965185210173401
$ws.Range("O18").Value = 1.002640018516937
$ws.Range("P18").Value = 0.9964334325644246

# Row 19: HexGrid-60degTilt5degRes
$ws.Range("A19").Value = 17
$ws.Range("A19").Font.Bold = $true
$ws.Range("A19").HorizontalAlignment = -4108
$ws.Range("A19").VerticalAlignment = -4160
$ws.Range("A19").Borders.LineStyle = 1
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C19").Value = 0.9410030045333887
$ws.Range("D19").Value = 1.114073650171882
$ws.Range("E19").Value = 0.971963238905884
$ws.Range("F19").Value = 1.00197140221889
$ws.Range("G19").Value = 0.9410030045333887
$ws.Range("H19").Value = 1.114073650171882
$ws.Range("I19").Value = 0.9380316020044707
$ws.Range("J19").Value = 1.017333132782812
$ws.Range("K19").Value = 0.952845175639116
$ws.Range("L19").Value = 1.071454969913451
$ws.Range("M19").Value = 0.9409757761636619
$ws.Range("N19").Value = 1.043018444538883
$ws.Range("O19").Value = 1.007252823957511
$ws.Range("P19").Value = 1.001084522021237
